$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Timestamps for data rows
$timestamps = @(
    "2021-10-05 13:40:11.719192",
    "2021-10-05 13:40:11.719206",
    "2021-10-05 13:40:11.719210",
    "2021-10-05 13:40:11.719213",
    "2021-10-05 13:40:11.719217",
    "2021-10-05 13:40:11.719221",
    "2021-10-05 13:40:11.719224",
    "2021-10-05 13:40:11.719228",
    "2021-10-05 13:40:11.719231",
    "2021-10-05 13:40:11.719234",
    "2021-10-05 13:40:11.719237",
    "2021-10-05 13:40:11.719241",
    "2021-10-05 13:40:11.719244",
    "2021-10-05 13:40:11.719247",
    "2021-10-05 13:40:11.719250"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
